# Weekly price-list refresh: a new observation (dated 2022-08-24 / serial
# 44797) is inserted as the new row 17 ("Bruselas (repollito)" at Mercado
# Mayorista Lo Valledor de Santiago), pushing the previously-existing rows
# 17-40 down to 18-41 (dimension grows from A1:R40 to A1:R41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 17, shifting rows 17:40 down to 18:41.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17, 1).Value  = 6
$ws.Cells.Item(17, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(17, 3).Value  = "Metropolitana"
$ws.Cells.Item(17, 4).Value  = 44797
$ws.Cells.Item(17, 5).Value  = 13
$ws.Cells.Item(17, 6).Value  = 100112035
$ws.Cells.Item(17, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(17, 8).Value  = "Sin especificar"
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 310
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = 18968
$ws.Cells.Item(17, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 1265
$ws.Cells.Item(17, 17).Value = 15
$ws.Cells.Item(17, 18).Value = "Hortaliza"
